$d = $word.ActiveDocument

# Helper: rewrite a paragraph's text as a single clean run, which also
# clears away any stray proofErr (spell-check) markers that were left
# over from splitting a word across runs.
function Fix-Paragraph($index, $newText) {
    $para = $d.Paragraphs.Item($index)
    $para.Range.Delete()
    $shifted = $d.Paragraphs.Item($index)
    $shifted.Range.InsertParagraphBefore()
    $target = $d.Paragraphs.Item($index)
    $target.Range.InsertBefore($newText)
}

# Paragraph indices (stable across these fixes since each one preserves
# the total paragraph count), found by inspecting $d.Paragraphs text.
Fix-Paragraph 32 "X =numpy.matrix([[4,5,1,2],[1,0,3,5],[2,1,8,2]])"
Fix-Paragraph 30 "import numpy"
Fix-Paragraph 15 "proc iml;"
Fix-Paragraph 4 "Abou El Dahab, Georges Michel; Shaver, Nicole; Nelson, John; Tuttle, Nathan"

# Move the _GoBack bookmark from the "GitHub link to Case Study 2 work"
# heading (its old location) to the blank paragraph right after the
# subtitle (its new location). Word only keeps a single _GoBack bookmark,
# so adding it here automatically removes the old one.
$para3 = $d.Paragraphs.Item(3)
$d.Bookmarks.Add("_GoBack", $para3.Range)
